# Updated cryptos list (Price / Volume(1h) columns), refreshed data pull.
# Rows 13/14 (WrappedEther <-> Polkadot) also swap rank order.
# Price cells that look like plain numbers are forced to Text format ("@")
# before assignment so Excel keeps the exact printed representation
# (e.g. "1.0000", "0.9993") instead of silently normalizing them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.197.43'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.847.93'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.86'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7030'
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07735'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3067'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.61'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07815'
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '93.29'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.140'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.851.27'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6868'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.592'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008333'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = '29.192.61'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.05'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '2.092.83'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.73'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1510'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.22'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.847'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.539'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.176'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.196'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05114'
$ws.Range('E33').Value = '  -1.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7895'
$ws.Range('E34').Value = '  +4.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.896'
$ws.Range('E35').Value = '  +2.95%  '
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.695'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').Value = '1.320.15'
$ws.Range('E38').Value = '  +7.86%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9588'
$ws.Range('E41').Value = '  +6.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.066'
$ws.Range('E42').Value = '  +9.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '106.89'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.693'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').Value = '1.992.68'
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5181'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.42'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.981'
$ws.Range('E51').Value = '  -0.74%  '

Write-Host "Applied crypto price updates"